# Load Screen code Updated
#
# 1) Cell C8 ("Run Mode" for Loads_TC001) changes from "Yes" to "NO".
# 2) A new test case row (row 10, Loads_TC003) is appended to the grid.
# 3) Selection / top-left-cell move down to the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update Run Mode for Loads_TC001 -------------------------------
$ws.Range("C8").Value = "NO"

# --- 2) Append new row 10: Loads_TC003 ---------------------------------
$newDescription = "Validate whehter carrier user can add Scoular loads for payment using full submit`n1) Enter valid user id and Password and click Login button.`n2) Click on Add New Load button from Load menu.`n3) Enter valid details in all required field and click Save button`n4) Now loads are saved successfully.`n5) Upload an Origin and Destination ticket image or PDF document for corresponding load.`n6) Observe Ready to Submit Load icon in grid should change to green color.`n7) Click on Submit Load button.`n8) Select any option and click Submit button."

$ws.Range("A10").Value = "Loads_TC003"
$ws.Range("A10").VerticalAlignment = -4108

$ws.Range("C10").Value = "YES"
$ws.Range("C10").VerticalAlignment = -4108

$ws.Range("B10").Value = $newDescription
$ws.Range("B10").WrapText = $true

$ws.Range("D10").Value = "Scoular loads validated successfully"
$ws.Range("D10").WrapText = $true
$ws.Range("D10").VerticalAlignment = -4108

$ws.Rows.Item(10).RowHeight = 180

# --- 3) Update the view so the new row is visible/selected ------------
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("D10").Select()
